$wb = $excel.ActiveWorkbook

# --- Rename shared strings used on the Sorting sheet ---
$wsSorting = $wb.Worksheets.Item("Sorting")
$wsSorting.Range("B3").Value = "Sauce Labs Onesie"
$wsSorting.Range("B5").Value = "Sauce Labs Backpack"

# --- Add a placeholder sheet at the end so the new copy gets a fresh sheetId ---
$placeholder = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# --- Duplicate the UserList sheet, placed after the placeholder (i.e. at the end) ---
$wsUserList = $wb.Worksheets.Item("UserList")
$wsUserList.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsUserList2 = $wb.Worksheets.Item($wb.Worksheets.Count)

# --- Remove the "locked_out_user" row (row 3) from the new sheet ---
$wsUserList2.Rows.Item(3).Delete()

# --- Remove the placeholder sheet now that the copy has taken the next sheetId ---
$placeholder.Delete()

# --- Re-fetch sheet references fresh (older handles can go stale after the delete) ---
$wsUserList = $wb.Worksheets.Item("UserList")
$wsSorting = $wb.Worksheets.Item("Sorting")
$wsUserList2 = $wb.Worksheets.Item("UserList (2)")

# --- Set the selection on the new sheet to match the recorded UI state ---
$wsUserList2.Activate()
$wsUserList2.Range("A3:XFD3").Select()

# --- Set the selection on the UserList sheet ---
$wsUserList.Activate()
$wsUserList.Range("A17").Select()

# --- Set the selection on the Sorting sheet, and leave it as the active sheet ---
$wsSorting.Activate()
$wsSorting.Range("D15").Select()
